$wb = $excel.ActiveWorkbook

$ws2 = $wb.Worksheets.Item("Tools_QA_Practice")

# Update continent value from "North America" to "South America"
$ws2.Range("J2").Value = "South America"

# Clear the Status value in L2 (was "Success"), keeping the cell entry present but empty
$ws2.Range("L2").ClearContents()
$ws2.Range("L2").Style = "Normal"

# Update the active selection on the Tools_QA_Practice sheet to K9
$ws2.Activate()
$ws2.Range("K9").Select()
